$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 9944
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 316.81
$ws.Range("D3").Value = 318.58999999999997
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = 0.56000000000000005
$ws.Range("G3").Value = 42606.427835648145
$ws.Range("H3").Value = $false

# Row 4
$ws.Range("A4").Value = 9875.39
$ws.Range("B4").Value = 9944
$ws.Range("C4").Value = 316.81
$ws.Range("D4").Value = 319
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 0.69
$ws.Range("G4").Value = 42606.48673611111
$ws.Range("H4").Value = $false
